$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1409.2727
$ws.Range("I40").Value = 1271.5714
$ws.Range("J40").Value = 1473.5333
$ws.Range("K40").Value = 1271.5714
$ws.Range("L40").Value = 1473.5333
$ws.Range("M40").Value = -1096.5714
$ws.Range("N40").Value = -1823.5333
$ws.Range("H98").Value = 1429.8572
$ws.Range("I98").Value = 1462.9231
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1462.9231
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 35.07690000000002
$ws.Range("N98").Value = -3996
$ws.Range("H107").Value = 5004.0557
$ws.Range("I107").Value = 6494.077
$ws.Range("K107").Value = 6494.077
$ws.Range("M107").Value = -4574.077
$ws.Range("H122").Value = 1429.8572
$ws.Range("I122").Value = 1462.9231
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4388.7693
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1938.7693
$ws.Range("N122").Value = -7900
$ws.Range("H139").Value = 30755.555
$ws.Range("J139").Value = 30755.555
$ws.Range("L139").Value = 30755.555
$ws.Range("N139").Value = -41035.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9248.039
$ws.Range("I32").Value = 9731.289
$ws.Range("J32").Value = 5623.6665
$ws.Range("K32").Value = 9731.289
$ws.Range("L32").Value = 5623.6665
$ws.Range("M32").Value = -9444.289
$ws.Range("N32").Value = -6197.6665
$ws.Range("H61").Value = 26319778
$ws.Range("I61").Value = 45458840
$ws.Range("K61").Value = 45458840
$ws.Range("M61").Value = -45458628
$ws.Range("H74").Value = 16132377
$ws.Range("I74").Value = 25001900
$ws.Range("K74").Value = 25001900
$ws.Range("M74").Value = -25001026
$ws.Range("H77").Value = 16132377
$ws.Range("I77").Value = 25001900
$ws.Range("K77").Value = 125009500
$ws.Range("M77").Value = -125005132
$ws.Range("H132").Value = 20836608
$ws.Range("I132").Value = 35716492
$ws.Range("K132").Value = 107149476
$ws.Range("M132").Value = -107146946
$ws.Range("H136").Value = 26319778
$ws.Range("I136").Value = 45458840
$ws.Range("K136").Value = 136376520
$ws.Range("M136").Value = -136373970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 782.7143
$ws.Range("I99").Value = 695.8
$ws.Range("K99").Value = 695.8
$ws.Range("M99").Value = 802.2
$ws.Range("H107").Value = 1332.1714
$ws.Range("I107").Value = 1339
$ws.Range("K107").Value = 1339
$ws.Range("M107").Value = 581
$ws.Range("H134").Value = 3411.2415
$ws.Range("I134").Value = 1846.1177
$ws.Range("K134").Value = 5538.3531
$ws.Range("M134").Value = -3003.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2291.25
$ws.Range("I105").Value = 1143.3334
$ws.Range("J105").Value = 2980
$ws.Range("K105").Value = 1143.3334
$ws.Range("L105").Value = 2980
$ws.Range("M105").Value = 603.6666
$ws.Range("N105").Value = -6474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 824.6875
$ws.Range("I97").Value = 803.7083
$ws.Range("J97").Value = 887.625
$ws.Range("K97").Value = 803.7083
$ws.Range("L97").Value = 887.625
$ws.Range("M97").Value = -307.7083
$ws.Range("N97").Value = -1879.625
$ws.Range("H102").Value = 2517.25
$ws.Range("I102").Value = 2874.0454
$ws.Range("J102").Value = 1209
$ws.Range("K102").Value = 2874.0454
$ws.Range("L102").Value = 1209
$ws.Range("M102").Value = -1252.0454
$ws.Range("N102").Value = -4453
$ws.Range("H113").Value = 167768.67
$ws.Range("I113").Value = 167768.67
$ws.Range("K113").Value = 167768.67
$ws.Range("M113").Value = -165598.67
$ws.Range("H126").Value = 3885.963
$ws.Range("J126").Value = 5093.2
$ws.Range("L126").Value = 15279.6
$ws.Range("N126").Value = -20219.6
$ws.Range("H132").Value = 6140.9443
$ws.Range("I132").Value = 5156.7144
$ws.Range("J132").Value = 6767.273
$ws.Range("K132").Value = 15470.1432
$ws.Range("L132").Value = 20301.819
$ws.Range("M132").Value = -12940.1432
$ws.Range("N132").Value = -25361.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2296.652
$ws.Range("I16").Value = 2296.652
$ws.Range("K16").Value = 2296.652
$ws.Range("M16").Value = -2126.652
$ws.Range("H40").Value = 5161.857
$ws.Range("I40").Value = 6839.3
$ws.Range("K40").Value = 6839.3
$ws.Range("M40").Value = -6703.3
$ws.Range("H82").Value = 1985.35
$ws.Range("I82").Value = 1744.8572
$ws.Range("J82").Value = 2546.5
$ws.Range("K82").Value = 1744.8572
$ws.Range("L82").Value = 2546.5
$ws.Range("M82").Value = -1383.8572
$ws.Range("N82").Value = -3268.5
$ws.Range("H85").Value = 1985.35
$ws.Range("I85").Value = 1744.8572
$ws.Range("J85").Value = 2546.5
$ws.Range("K85").Value = 1744.8572
$ws.Range("L85").Value = 2546.5
$ws.Range("M85").Value = -496.8571999999999
$ws.Range("N85").Value = -5042.5
$ws.Range("H132").Value = 10211226
$ws.Range("I132").Value = 4950.591
$ws.Range("J132").Value = 18527452
$ws.Range("K132").Value = 14851.773
$ws.Range("L132").Value = 55582356
$ws.Range("M132").Value = -12321.773
$ws.Range("N132").Value = -55587416
$ws.Range("H136").Value = 15157906
$ws.Range("I136").Value = 20835230
$ws.Range("J136").Value = 18373.334
$ws.Range("K136").Value = 62505690
$ws.Range("L136").Value = 55120.00199999999
$ws.Range("M136").Value = -62503140
$ws.Range("N136").Value = -60220.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1489.9333
$ws.Range("I113").Value = 1007.5
$ws.Range("J113").Value = 2213.5833
$ws.Range("K113").Value = 3022.5
$ws.Range("L113").Value = 6640.749899999999
$ws.Range("M113").Value = -852.5
$ws.Range("N113").Value = -10980.7499
$ws.Range("H122").Value = 3046.8235
$ws.Range("I122").Value = 2929.3572
$ws.Range("K122").Value = 8788.0716
$ws.Range("M122").Value = -6338.071599999999
$ws.Range("H132").Value = 1941.5
$ws.Range("I132").Value = 973.7273
$ws.Range("K132").Value = 2921.1819
$ws.Range("M132").Value = -391.1819
$ws.Range("H136").Value = 799.2121
$ws.Range("I136").Value = 777.3125
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2331.9375
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 218.0625
$ws.Range("N136").Value = -9600
